$wb = $excel.ActiveWorkbook

# OFF sheet - row 2 (H) updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 385
$wsOff.Range("C2").Value = 270
$wsOff.Range("D2").Value = 86
$wsOff.Range("E2").Value = 39

# DEF sheet - row 2 (H) updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 369
$wsDef.Range("C2").Value = 262
$wsDef.Range("D2").Value = 78
$wsDef.Range("E2").Value = 26
